# Updated symbol list on Mon Dec 12 17:31:13 UTC 2022 with GitHub Actions
#
# Column D values are stored as literal text (not numbers) in the workbook,
# so numeric-looking updates are entered with a leading apostrophe (forces
# text entry, like a user typing '276.38 into Excel) and then ClearFormats()
# is used to drop the transient "Text" / quote-prefix cell style that Excel
# applies to quote-prefixed entries, so the cell's style index is left
# unchanged (matching the source diff, which touches no styles).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'276.38"
$ws.Range("D2").ClearFormats()

$ws.Range("D3").Value = "'20.99"
$ws.Range("D3").ClearFormats()

$ws.Range("D4").Value = "'6.216"
$ws.Range("D4").ClearFormats()

$ws.Range("D5").Value = "'0.06186"
$ws.Range("D5").ClearFormats()

$ws.Range("D7").Value = "'1.517"
$ws.Range("D7").ClearFormats()

$ws.Range("D8").Value = "'6.548"
$ws.Range("D8").ClearFormats()

$ws.Range("D10").Value = "'0.1639"
$ws.Range("D10").ClearFormats()

$ws.Range("D11").Value = "'0.08207"
$ws.Range("D11").ClearFormats()

$ws.Range("D12").Value = "'0.03422"
$ws.Range("D12").ClearFormats()

$ws.Range("D13").Value = "'0.03129"
$ws.Range("D13").ClearFormats()

$ws.Range("D14").Value = "'0.09130"
$ws.Range("D14").ClearFormats()

$ws.Range("D15").Value = "'3.769"
$ws.Range("D15").ClearFormats()

$ws.Range("D16").Value = "'0.001609"
$ws.Range("D16").ClearFormats()

$ws.Range("D18").Value = "'0.006275"
$ws.Range("D18").ClearFormats()

$ws.Range("D19").Value = "'0.006139"
$ws.Range("D19").ClearFormats()

$ws.Range("D21").Value = "'0.0001500"
$ws.Range("D21").ClearFormats()

$ws.Range("D22").Value = "'3.738"
$ws.Range("D22").ClearFormats()

$ws.Range("D25").Value = "'0.3285"
$ws.Range("D25").ClearFormats()

$ws.Range("D28").Value = "'0.0002739"
$ws.Range("D28").ClearFormats()

$ws.Range("D40").Value = "'0.04670"
$ws.Range("D40").ClearFormats()

$ws.Range("D41").Value = "'0.007020"
$ws.Range("D41").ClearFormats()

# Row 42 and row 43 swap coins (BKEXToken <-> CEJI), each keeping its own
# rank prefix in column E ("41..." stays on row 42, "42..." stays on row 43).
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003521"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1104"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "42BKEXTokenBKK"

$ws.Range("D44").Value = "'0.01114"
$ws.Range("D44").ClearFormats()

$ws.Range("D45").Value = "'0.00006418"
$ws.Range("D45").ClearFormats()

$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").ClearFormats()

$ws.Range("D47").Value = "'0.8459"
$ws.Range("D47").ClearFormats()

$ws.Range("D48").Value = "'0.001385"
$ws.Range("D48").ClearFormats()

$ws.Range("D49").Value = "'0.00001902"
$ws.Range("D49").ClearFormats()
